# Sync the "Price" (D) / "Volume(1h)" (E) / "Coin" (B) / "Link" (C) columns
# of the crypto ranking sheet to the latest scraped snapshot.
#
# The source data stores these as literal text (e.g. "0.06940", "-0.80%")
# rather than numbers, to preserve exact formatting (trailing zeros, the
# percent sign, etc.). Plain numeric-looking strings assigned through
# Range.Value are auto-converted to numbers by Excel, so for the
# numeric-looking cells (Price / Volume) we first set NumberFormat to
# Text ("@") to force the literal text to be preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.22"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.78%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.27"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.86%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.118"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.63%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05686"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.52%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.517"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.63%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8195"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.83%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8576"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.40%"
$ws.Range("B9").Value = "MandalaExchangeToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06940"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.80%"
$ws.Range("B10").Value = "BitrueCoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.02857"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.69%"
$ws.Range("B11").Value = "BitMartToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09393"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.09%"
$ws.Range("B12").Value = "BitForexToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.001531"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.77%"
$ws.Range("B13").Value = "CoinExToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04009"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-13.90%"
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0005989"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-93.97%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006215"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.03%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.512"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.58%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.009"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.33%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.230"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "8.52%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3149"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.17%"
$ws.Range("B20").Value = "WazirX"
$ws.Range("C20").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1333"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.43%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.03223"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.03%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.04%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.565"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-4.75%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.01%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001216"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.15%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004464"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-2.45%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001179"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "22.85%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-27.43%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03717"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.50%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1058"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.07%"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.002419"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.25%"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.005977"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.81%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009709"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "17.42%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005110"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-5.26%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.01%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-30.35%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "3.23%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.01%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.01%"
